# Edit script implementing the diff:
#  - Delete the "Ghost" bullet paragraph entirely.
#  - Reorder/rewrite the "Nirvana" and "Reincarnation" bullets: the first
#    bullet becomes "Reincarnation" (merging in the old Reincarnation body
#    text), and a new second bullet "Nirvana" follows with new transition
#    wording leading into the original Nirvana payoff text.
#  - Rewrite the start of the "Purgatory" bullet and move the _GoBack
#    bookmark there (removing it from the end of "Torment").
#  - Move the lastRenderedPageBreak rendering hint from the "player chooses
#    yes" paragraph to the "player chooses no" paragraph.

$d = $word.ActiveDocument
$dash = [char]0x2013
$apos = [char]0x2019

# ---------------------------------------------------------------------
# Step 1: delete the entire "Ghost" bullet paragraph (paragraph 15).
# ---------------------------------------------------------------------
$ghostPara = $d.Paragraphs.Item(15)
$ghostPara.Range.Delete()

# ---------------------------------------------------------------------
# Step 2: split paragraph 13 (currently the "Nirvana" bullet) right before
# "What if the afterlife was a sort of Nirvana?" so the first half keeps
# the intro sentences and the second half keeps the Nirvana payoff text.
# ---------------------------------------------------------------------
$p13 = $d.Paragraphs.Item(13)
$full13 = $p13.Range.Text
$marker = "in store?  "
$splitIdx = $full13.IndexOf($marker) + $marker.Length
$splitPoint = $p13.Range.Start + $splitIdx
$splitRange = $d.Range($splitPoint, $splitPoint)
$splitRange.InsertParagraphAfter() | Out-Null

# Now:
#   13 = "Nirvana - This paradise is what most seem ... in store?  "
#   14 = "What if the afterlife was a sort of Nirvana?  What if we become ... goal."
#   15 = "Reincarnation - Perhaps there isn't even an afterlife. ... instruments."

# ---------------------------------------------------------------------
# Step 3: turn paragraph 13 into the new "Reincarnation" bullet: change
# the bold heading text, then append the merged-in Reincarnation body.
# ---------------------------------------------------------------------
$p13 = $d.Paragraphs.Item(13)
$headingRange = $p13.Range.Duplicate
$headingRange.Find.Execute("Nirvana " + $dash + " ", $true, $false, $false, $false, $false, $true, 1, $false, "Reincarnation " + $dash + " ", 2) | Out-Null

$p13 = $d.Paragraphs.Item(13)
$appendPoint = $p13.Range.End - 1
$appendRange = $d.Range($appendPoint, $appendPoint)
$appendText = "Perhaps there isn" + $apos + "t even an afterlife.  What if after death, our souls were reincarnated?  Maybe our world is full of souls thousands of years old; the same people in different bodies, like the same melody played on different instruments."
$appendRange.InsertAfter($appendText) | Out-Null

# ---------------------------------------------------------------------
# Step 4: turn paragraph 14 into the new "Nirvana" bullet: prepend a bold
# "Nirvana -" heading run, and swap the lead-in wording.
# ---------------------------------------------------------------------
$p14 = $d.Paragraphs.Item(14)
$leadRange = $p14.Range.Duplicate
$leadRange.Find.Execute("What if the afterlife was", $true, $false, $false, $false, $false, $true, 1, $false, "Perhaps, the afterlife is", 2) | Out-Null

$p14 = $d.Paragraphs.Item(14)
$insertPoint = $d.Range($p14.Range.Start, $p14.Range.Start)
$insertPoint.InsertAfter("Nirvana " + $dash) | Out-Null
$headingOnly = $d.Range($p14.Range.Start, $p14.Range.Start + ("Nirvana " + $dash).Length)
$headingOnly.Bold = 1

Write-Host "13:" $d.Paragraphs.Item(13).Range.Text
Write-Host "14:" $d.Paragraphs.Item(14).Range.Text
Write-Host "15:" $d.Paragraphs.Item(15).Range.Text
